# Accessibility Chart.xlsx - apply commit changes via Excel COM-interop
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1 edits ---

# Row 4 (Tech column total): was 5, now 4
$ws1.Range("L4").Value = 4

# Rows 23:43 - the "Tech" column (L) gains an "X" mark matching column K,
# for every row that already has an X in K (same value + same style/format).
$ws1.Range("K23:K43").Copy($ws1.Range("L23:L43"))

# Update Sheet1's saved view/selection (no longer scrolled to A46 / B36 selected)
$ws1.Range("K7").Select()

# --- Sheet2 edits ---

# Row 2: "Employed" sub-item renamed from "ID #"/"employed" to "UserId"/"uID"
$ws2.Range("C2").Value = "UserId"
$ws2.Range("G2").Value = "uID"

# Update Sheet2's saved selection, and re-activate Sheet2 so it remains the
# active tab in the saved workbook (selecting on Sheet1 above made it active).
$ws2.Range("G3").Select()
